# "added fresh data for registration"
# Refresh the registration test data on the "RegDetails" sheet with a new
# batch of emails / names, then move the active selection to E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New email addresses (column D / EmailAddress) for rows 2-5.
$ws.Range("D2").Value = "a25@email.com"
$ws.Range("D3").Value = "a26@email.com"
$ws.Range("D4").Value = "a27@email.com"
$ws.Range("D5").Value = "a28@email.com"

# New middle name / last name pairs (columns B and C) for rows 2-5.
$ws.Range("B2").Value = "Ronald"
$ws.Range("C2").Value = "Delver"

$ws.Range("B3").Value = "Peter"
$ws.Range("C3").Value = "Con"

$ws.Range("B4").Value = "Asult"
$ws.Range("C4").Value = "Bolswa"

$ws.Range("B5").Value = "Jimmy"
$ws.Range("C5").Value = "Lever"

# Resize the book view the way Excel records it after the edit.
$excel.ActiveWindow.Width = 19380
$excel.ActiveWindow.Height = 5955

# Move the active selection to E4, matching the post-edit cursor position.
$ws.Range("E4").Select()
